$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (Coin price column D) - must remain text, not auto-converted to numbers
$textUpdates = @(
    @('D2', '98.842.54'),
    @('D3', '3.373.66'),
    @('D5', '259.45'),
    @('D6', '628.85'),
    @('D8', '0.393'),
    @('D10', '0.865'),
    @('D11', '3.372.82'),
    @('D13', '98.614.38'),
    @('D14', '36.27'),
    @('D15', '0.0000249'),
    @('D16', '4.000.45'),
    @('D17', '5.49'),
    @('D18', '3.373.26'),
    @('D19', '3.56'),
    @('D20', '15.25'),
    @('D21', '490.53'),
    @('D23', '0.0000212'),
    @('D24', '9.37'),
    @('D25', '5.65'),
    @('D26', '88.81'),
    @('D28', '3.554.13'),
    @('D29', '0.285'),
    @('D31', '0.191'),
    @('D32', '0.135'),
    @('D33', '9.69'),
    @('D34', '0.998'),
    @('D35', '28.06'),
    @('D36', '0.151'),
    @('D37', '7.32'),
    @('D39', '500.31'),
    @('D40', '0.461'),
    @('D41', '24.91'),
    @('D42', '3.76'),
    @('D43', '1.27'),
    @('D44', '3.31'),
    @('D45', '0.788'),
    @('D47', '160.09'),
    @('D48', '1.95'),
    @('D49', '0.841'),
    @('D50', '4.64')
)

foreach ($item in $textUpdates) {
    $cellRef = $item[0]
    $val = $item[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Plain cell updates (Coin name, Link, Volume columns) - Excel keeps these as text naturally
$plainUpdates = @(
    @('E2', '  +1.51%  '),
    @('E3', '  +7.44%  '),
    @('E4', '  +0.09%  '),
    @('E5', '  +6.91%  '),
    @('E6', '  +2.61%  '),
    @('E7', '  +24.07%  '),
    @('E8', '  +1.95%  '),
    @('E9', '  +0.01%  '),
    @('E10', '  +10.39%  '),
    @('E11', '  +7.51%  '),
    @('E12', '  +0.12%  '),
    @('E13', '  +1.63%  '),
    @('E14', '  +6.17%  '),
    @('E15', '  +3.28%  '),
    @('E16', '  +7.34%  '),
    @('E17', '  -0.11%  '),
    @('E18', '  +7.28%  '),
    @('E19', '  -0.19%  '),
    @('E20', '  +4.04%  '),
    @('E21', '  -6.03%  '),
    @('E22', '  +7.74%  '),
    @('E23', '  +9.17%  '),
    @('E24', '  +6.12%  '),
    @('E25', '  +2.56%  '),
    @('E26', '  -0.07%  '),
    @('E27', '  +2.29%  '),
    @('E28', '  +7.59%  '),
    @('E29', '  +18.22%  '),
    @('E30', '  -0.10%  '),
    @('E31', '  +8.83%  '),
    @('E32', '  +8.92%  '),
    @('B33', 'InternetComputer(DFINITY)'),
    @('C33', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('E33', '  +7.17%  '),
    @('B34', 'Binance-PegBSC-USD'),
    @('C34', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'),
    @('E34', '  +7.47%  '),
    @('E35', '  +4.93%  '),
    @('B36', 'Kaspa'),
    @('C36', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'),
    @('E36', '  -1.50%  '),
    @('B37', 'RenderToken'),
    @('C37', 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'),
    @('E37', '  -1.39%  '),
    @('E38', '  +3.84%  '),
    @('E39', '  +4.55%  '),
    @('E40', '  +5.12%  '),
    @('E41', '  +2.71%  '),
    @('E42', '  +3.73%  '),
    @('E43', '  +3.08%  '),
    @('E44', '  +4.33%  '),
    @('E46', '  +0.03%  '),
    @('E47', '  -0.38%  '),
    @('E48', '  +0.56%  '),
    @('E49', '  +12.79%  '),
    @('E50', '  +2.86%  '),
    @('E51', '  +4.25%  ')
)

foreach ($item in $plainUpdates) {
    $ws.Range($item[0]).Value = $item[1]
}
